# Apply the changes described by the diff to the Tab11 worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab11")

# 1) Fix mis-encoded accented characters in the "Regional Economic Communities" footnote (cell A103).
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# 2) Small recalculation tweak to I67.
$ws.Range("I67").Value = 62.926791416041397

# 3) Updated aggregate figures for row 97 ("Africa, Fragile States").
$ws.Range("C97").Value = 0.56318884957869997
$ws.Range("D97").Value = 26.537375356280499
$ws.Range("E97").Value = 76.474112224555299
$ws.Range("F97").Value = 0.73576110211611001
$ws.Range("G97").Value = 28.588177796322299
$ws.Range("H97").Value = 842905.22956505604
$ws.Range("I97").Value = 57.498311168908401
$ws.Range("J97").Value = 25.267971338376999

# 4) Updated aggregate figures for row 98 ("ROW, Fragile States").
$ws.Range("C98").Value = 4.1590685805839298
$ws.Range("D98").Value = 32.135292849530899
$ws.Range("E98").Value = 95.525715488116603
$ws.Range("F98").Value = 6.8758422491245801
$ws.Range("G98").Value = 50.900082808909602
$ws.Range("H98").Value = 1792876.9873390901
$ws.Range("I98").Value = 80.413810560134095
$ws.Range("J98").Value = 29.6264260164983

$wb.Save()
